$d = $word.ActiveDocument

$replacements = @(
    @("687×2=", "714×2="),
    @("279×7=", "958×2="),
    @("661×4=", "797×4="),
    @("280×4=", "711×5="),
    @("400×9=", "914×3="),
    @("402×9=", "979×5="),
    @("497×7=", "938×7="),
    @("695×7=", "309×2="),
    @("431×7=", "565×6="),
    @("525×8=", "666×3="),
    @("142×4=", "287×8="),
    @("961×4=", "809×6="),
    @("842×5=", "706×6="),
    @("769×2=", "685×6="),
    @("178×4=", "908×8="),
    @("694×5=", "988×6="),
    @("684×7=", "917×8="),
    @("300×5=", "268×6="),
    @("698×2=", "931×2="),
    @("650×4=", "903×2="),
    @("621×2=", "667×6="),
    @("632×7=", "973×7="),
    @("875×4=", "934×2="),
    @("488×4=", "409×4="),
    @("600×8=", "716×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
